$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.259725451469421
$ws.Range("B1").Value = 1.973381042480469
$ws.Range("C1").Value = 5.757643222808838
$ws.Range("D1").Value = 1.951203346252441
$ws.Range("E1").Value = 1.126999974250793
